$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header row styles: P1 reverts to style 1 (like A1), Q1 reverts to style 2 (like O1) ---
$ws.Range("A1").Copy()
$ws.Range("P1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("O1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 2 ---
$ws.Range("A2").Value = 484848484
$ws.Range("B2").Value = "Dobbel 15 GB"
$ws.Range("C2").Value = "Fornavn"
$ws.Range("D2").Value = "Etternavn"
$ws.Range("E2").Value = "Adresaa"
$ws.Range("F2").Value = "Adresa"
$ws.Range("G2").Value = "Adresa"
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = "1@1"
$ws.Range("M2").Value = "1@1"
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = "Faktura 1"

# --- Row 3 ---
$ws.Range("A3").Value = 487897214
$ws.Range("B3").Value = "Dobbel 15 GB"
$ws.Range("C3").Value = "Fornavn"
$ws.Range("D3").Value = "Etternavn"
$ws.Range("E3").Value = "Adresaa"
$ws.Range("F3").Value = "Adresa"
$ws.Range("G3").Value = "Adresa"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = "2@2"
$ws.Range("M3").Value = "2@2"
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = "Faktura 2"

# --- Row 4 (new row) ---
$ws.Range("A4").Value = 48487897
$ws.Range("B4").Value = "Dobbel 15 GB"
$ws.Range("C4").Value = "Fornavn"
$ws.Range("D4").Value = "Etternavn"
$ws.Range("E4").Value = "EHF"
$ws.Range("F4").Value = "Adresa"
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 4
$ws.Range("L4").Value = "3@3"
$ws.Range("M4").Value = "3@3"
$ws.Range("N4").Value = 2
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = "Faktura 3"

# --- Hyperlinks ---
# Remove old hyperlinks on L2/M2/L3/M3 so we can re-add with the new mail targets
$ws.Range("L2").Hyperlinks.Delete()
$ws.Range("M2").Hyperlinks.Delete()
$ws.Range("L3").Hyperlinks.Delete()
$ws.Range("M3").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("L2"), "mailto:1@1", [Type]::Missing, [Type]::Missing, "1@1")
$ws.Hyperlinks.Add($ws.Range("M2"), "mailto:1@1", [Type]::Missing, [Type]::Missing, "1@1")
$ws.Hyperlinks.Add($ws.Range("L3"), "mailto:2@2", [Type]::Missing, [Type]::Missing, "2@2")
$ws.Hyperlinks.Add($ws.Range("M3"), "mailto:2@2", [Type]::Missing, [Type]::Missing, "2@2")
$ws.Hyperlinks.Add($ws.Range("L4"), "mailto:3@3", [Type]::Missing, [Type]::Missing, "3@3")
$ws.Hyperlinks.Add($ws.Range("M4"), "mailto:3@3", [Type]::Missing, [Type]::Missing, "3@3")

# Re-apply the original hyperlink cell format (style index 3) since Hyperlinks.Add
# creates a near-duplicate style; copy format from a cell that already has it.
$ws.Range("L2").Copy()
$ws.Range("L2:M4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Sheet view: select E4, clear frozen/scrolled topLeftCell ---
$ws.Range("A1").Select() | Out-Null
$ws.Range("E4").Select() | Out-Null

$wb.Save()
